# Trade #74 closed at 2026-02-17 12:57:23 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.23
$summary.Range("B4").Value = 0.22
$summary.Range("B5").Value = 0.06
$summary.Range("B6").Value = 74
$summary.Range("B8").Value = 23
$summary.Range("B9").Value = 44.59

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.23
$status.Range("D4").Value = 74
$status.Range("E4").Value = 0.22
$status.Range("F4").Value = 0.23
$status.Range("G4").Value = 44.59

# --- Append new trade row (row 75) to both "All Trades" and "MarketMaking" sheets ---
$newRow = @{
    A = 74
    B = "2026-02-17"
    C = "12:57:16"
    D = "MarketMaking"
    E = "UP"
    F = 0.6
    G = 0.57
    H = "CLOSED"
    I = -5
    J = -0.03
    K = 100.23
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A75").Value = $newRow.A
    # Force text format so the date-like string "2026-02-17" is not
    # auto-converted into a date serial number.
    $ws.Range("B75").NumberFormat = "@"
    $ws.Range("B75").Value = $newRow.B
    $ws.Range("C75").Value = $newRow.C
    $ws.Range("D75").Value = $newRow.D
    $ws.Range("E75").Value = $newRow.E
    $ws.Range("F75").Value = $newRow.F
    $ws.Range("G75").Value = $newRow.G
    $ws.Range("H75").Value = $newRow.H
    $ws.Range("I75").Value = $newRow.I
    $ws.Range("J75").Value = $newRow.J
    $ws.Range("K75").Value = $newRow.K
    $ws.Range("L75").Value = $newRow.L
    $ws.Range("M75").Value = $newRow.M
    $ws.Range("N75").Value = $newRow.N
    $ws.Range("O75").Value = $newRow.O
    $ws.Range("P75").Value = $newRow.P
    $ws.Range("Q75").Value = $newRow.Q
}
